# 14/01 - adicionar melhor goleiro
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new names used only in the new rows (44-65) first, column A,
#        row by row, so the shared-string table grows in the same order the
#        original authoring session produced it in (new names before the new
#        "Melhor Goleiro" header string).
$newRowNames = @(
    "Romario",
    "Tom",
    "Ismael",
    "Cabeleira",
    "Athos",
    "Fabinho",
    "Leandrão",
    "Corinthiano",
    "Leandrinho",
    "Miqueias",
    "Michel",
    "Peixe",
    "Geovane",
    "Du",
    "Fernando",
    "Digão",
    "Marcelão",
    "Leah",
    "Eduardo",
    "Juscielio",
    "Milton",
    "Lucian"
)
for ($i = 0; $i -lt $newRowNames.Length; $i++) {
    $ws.Cells.Item(44 + $i, 1).Value = $newRowNames[$i]
}

# --- 2. New column L header ("Melhor Goleiro") - becomes the last new shared
#        string, matching the source workbook.
$ws.Range("L1").Value = "Melhor Goleiro"

# --- 3. Fill in column L (Melhor Goleiro) with 0 for every already-existing
#        player row (2-43).
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 12).Value = 0
}

# --- 4. Fill the remaining stats (C..L) for the new rows (44-65).
$rows = @(
    @(44, 2, 1, 4, 0, 1, 0, 0, 0, 0, 0),
    @(45, 2, 1, 4, 1, 1, 0, 0, 0, 0, 0),
    @(46, 2, 1, 4, 0, 1, 0, 0, 0, 0, 0),
    @(47, 2, 1, 4, 2, 1, 0, 0, 0, 0, 0),
    @(48, 2, 1, 4, 1, 1, 0, 0, 0, 0, 0),
    @(49, 0, 3, 4, 2, 1, 0, 1, 0, 0, 0),
    @(50, 0, 3, 4, 0, 1, 0, 1, 0, 0, 0),
    @(51, 0, 3, 4, 0, 1, 0, 1, 0, 0, 0),
    @(52, 0, 3, 4, 1, 1, 0, 1, 0, 0, 0),
    @(53, 0, 3, 4, 0, 1, 0, 1, 0, 0, 0),
    @(54, 2, 3, 2, 2, 1, 0, 0, 0, 0, 0),
    @(55, 2, 3, 2, 2, 1, 0, 0, 0, 0, 0),
    @(56, 2, 3, 2, 2, 1, 0, 0, 0, 0, 0),
    @(57, 2, 3, 2, 0, 1, 0, 0, 0, 0, 0),
    @(58, 2, 3, 2, 2, 1, 0, 0, 0, 0, 0),
    @(59, 6, 3, 0, 1, 1, 1, 0, 0, 0, 0),
    @(60, 6, 3, 0, 2, 1, 1, 0, 0, 0, 0),
    @(61, 6, 3, 0, 3, 1, 1, 0, 1, 0, 0),
    @(62, 6, 3, 0, 3, 1, 1, 0, 0, 0, 0),
    @(63, 6, 3, 0, 3, 1, 1, 0, 0, 0, 0),
    @(64, 6, 5, 3, 3, 1, 1, 0, 0, 12, 1),
    @(65, 3, 5, 6, 0, 1, 0, 1, 0, 18, 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value  = $row[1]   # C Vitorias
    $ws.Cells.Item($r, 4).Value  = $row[2]   # D Empate
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E Derrotas
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F Gols
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G Partidas
    $ws.Cells.Item($r, 8).Value  = $row[6]   # H Tarde de Vitoria
    $ws.Cells.Item($r, 9).Value  = $row[7]   # I La barca
    $ws.Cells.Item($r, 10).Value = $row[8]   # J Craque do Dia
    $ws.Cells.Item($r, 11).Value = $row[9]   # K Gols Sofridos
    $ws.Cells.Item($r, 12).Value = $row[10]  # L Melhor Goleiro
}

# --- 5. Drop the AutoFilter that covered the old range.
$ws.AutoFilterMode = $false

# --- 6. Restore the selection to the new header cell, as in the source file.
$ws.Range("L1").Select()
